$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 123
$ws.Range("A3").Value = 345
$ws.Range("A4").Value = 124
$ws.Range("A5").Value = 30697813
$ws.Range("B5").Value = 86.5

$ws.PageSetup.Orientation = 1

$ws.Range("B5").Select()
